$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.244075179100037
$ws.Range("B1").Value = 2.375225067138672
$ws.Range("C1").Value = 3.931496381759644
$ws.Range("D1").Value = 2.983071804046631
$ws.Range("E1").Value = 1.29889714717865
